$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "=""58.584.20"""
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Formula = "=""  -1.59%  """
$ws.Range("E2").Copy()
$ws.Range("E2").PasteSpecial(-4163)
$ws.Range("D3").Formula = "=""2.623.40"""
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Formula = "=""  +0.78%  """
$ws.Range("E3").Copy()
$ws.Range("E3").PasteSpecial(-4163)
$ws.Range("D4").Formula = "=""1.00"""
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Formula = "=""  +0.00%  """
$ws.Range("E4").Copy()
$ws.Range("E4").PasteSpecial(-4163)
$ws.Range("D5").Formula = "=""534.86"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Formula = "=""  -0.79%  """
$ws.Range("E5").Copy()
$ws.Range("E5").PasteSpecial(-4163)
$ws.Range("D6").Formula = "=""143.36"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Formula = "=""  +1.15%  """
$ws.Range("E6").Copy()
$ws.Range("E6").PasteSpecial(-4163)
$ws.Range("E7").Formula = "=""  -0.03%  """
$ws.Range("E7").Copy()
$ws.Range("E7").PasteSpecial(-4163)
$ws.Range("D8").Formula = "=""0.571"""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Formula = "=""  +0.92%  """
$ws.Range("E8").Copy()
$ws.Range("E8").PasteSpecial(-4163)
$ws.Range("E9").Formula = "=""  +5.44%  """
$ws.Range("E9").Copy()
$ws.Range("E9").PasteSpecial(-4163)
$ws.Range("E10").Formula = "=""  -2.44%  """
$ws.Range("E10").Copy()
$ws.Range("E10").PasteSpecial(-4163)
$ws.Range("D11").Formula = "=""0.332"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Formula = "=""  -0.74%  """
$ws.Range("E11").Copy()
$ws.Range("E11").PasteSpecial(-4163)
$ws.Range("E12").Formula = "=""  +1.33%  """
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial(-4163)
$ws.Range("D13").Formula = "=""3.086.03"""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Formula = "=""  +0.78%  """
$ws.Range("E13").Copy()
$ws.Range("E13").PasteSpecial(-4163)
$ws.Range("D14").Formula = "=""58.496.91"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Formula = "=""  -1.59%  """
$ws.Range("E14").Copy()
$ws.Range("E14").PasteSpecial(-4163)
$ws.Range("D15").Formula = "=""20.74"""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Formula = "=""  -0.32%  """
$ws.Range("E15").Copy()
$ws.Range("E15").PasteSpecial(-4163)
$ws.Range("D16").Formula = "=""2.630.96"""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Formula = "=""  +1.74%  """
$ws.Range("E16").Copy()
$ws.Range("E16").PasteSpecial(-4163)
$ws.Range("E17").Formula = "=""  -1.44%  """
$ws.Range("E17").Copy()
$ws.Range("E17").PasteSpecial(-4163)
$ws.Range("D19").Formula = "=""334.56"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Formula = "=""  -1.96%  """
$ws.Range("E19").Copy()
$ws.Range("E19").PasteSpecial(-4163)
$ws.Range("D20").Formula = "=""10.16"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Formula = "=""  +0.54%  """
$ws.Range("E20").Copy()
$ws.Range("E20").PasteSpecial(-4163)
$ws.Range("E21").Formula = "=""  -1.73%  """
$ws.Range("E21").Copy()
$ws.Range("E21").PasteSpecial(-4163)
$ws.Range("E22").Formula = "=""  +0.15%  """
$ws.Range("E22").Copy()
$ws.Range("E22").PasteSpecial(-4163)
$ws.Range("D23").Formula = "=""66.32"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Formula = "=""  -1.45%  """
$ws.Range("E23").Copy()
$ws.Range("E23").PasteSpecial(-4163)
$ws.Range("D24").Formula = "=""0.418"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Formula = "=""  +2.23%  """
$ws.Range("E24").Copy()
$ws.Range("E24").PasteSpecial(-4163)
$ws.Range("E25").Formula = "=""  -1.37%  """
$ws.Range("E25").Copy()
$ws.Range("E25").PasteSpecial(-4163)
$ws.Range("D26").Formula = "=""0.998"""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Formula = "=""  +0.10%  """
$ws.Range("E26").Copy()
$ws.Range("E26").PasteSpecial(-4163)
$ws.Range("D27").Formula = "=""7.12"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Formula = "=""  -1.42%  """
$ws.Range("E27").Copy()
$ws.Range("E27").PasteSpecial(-4163)
$ws.Range("D28").Formula = "=""0.0₃0736"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Formula = "=""  -1.26%  """
$ws.Range("E28").Copy()
$ws.Range("E28").PasteSpecial(-4163)
$ws.Range("D29").Formula = "=""0.999"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Formula = "=""  +0.00%  """
$ws.Range("E29").Copy()
$ws.Range("E29").PasteSpecial(-4163)
$ws.Range("E30").Formula = "=""  -2.15%  """
$ws.Range("E30").Copy()
$ws.Range("E30").PasteSpecial(-4163)
$ws.Range("D31").Formula = "=""5.90"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Formula = "=""  +0.49%  """
$ws.Range("E31").Copy()
$ws.Range("E31").PasteSpecial(-4163)
$ws.Range("D32").Formula = "=""18.91"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Formula = "=""  +0.43%  """
$ws.Range("E32").Copy()
$ws.Range("E32").PasteSpecial(-4163)
$ws.Range("D33").Formula = "=""151.28"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Formula = "=""  +0.90%  """
$ws.Range("E33").Copy()
$ws.Range("E33").PasteSpecial(-4163)
$ws.Range("E34").Formula = "=""  -2.18%  """
$ws.Range("E34").Copy()
$ws.Range("E34").PasteSpecial(-4163)
$ws.Range("D35").Formula = "=""0.854"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Formula = "=""  +0.74%  """
$ws.Range("E35").Copy()
$ws.Range("E35").PasteSpecial(-4163)
$ws.Range("E36").Formula = "=""  -2.00%  """
$ws.Range("E36").Copy()
$ws.Range("E36").PasteSpecial(-4163)
$ws.Range("E37").Formula = "=""  -1.83%  """
$ws.Range("E37").Copy()
$ws.Range("E37").PasteSpecial(-4163)
$ws.Range("E38").Formula = "=""  -2.96%  """
$ws.Range("E38").Copy()
$ws.Range("E38").PasteSpecial(-4163)
$ws.Range("E39").Formula = "=""  +0.85%  """
$ws.Range("E39").Copy()
$ws.Range("E39").PasteSpecial(-4163)
$ws.Range("D40").Formula = "=""282.12"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Formula = "=""  +3.25%  """
$ws.Range("E40").Copy()
$ws.Range("E40").PasteSpecial(-4163)
$ws.Range("D41").Formula = "=""1.00"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Formula = "=""  +0.01%  """
$ws.Range("E41").Copy()
$ws.Range("E41").PasteSpecial(-4163)
$ws.Range("D42").Formula = "=""0.595"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Formula = "=""  -1.03%  """
$ws.Range("E42").Copy()
$ws.Range("E42").PasteSpecial(-4163)
$ws.Range("E43").Formula = "=""  -0.48%  """
$ws.Range("E43").Copy()
$ws.Range("E43").PasteSpecial(-4163)
$ws.Range("D44").Formula = "=""19.05"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Formula = "=""  +2.45%  """
$ws.Range("E44").Copy()
$ws.Range("E44").PasteSpecial(-4163)
$ws.Range("E45").Formula = "=""  -1.41%  """
$ws.Range("E45").Copy()
$ws.Range("E45").PasteSpecial(-4163)
$ws.Range("E46").Formula = "=""  +0.63%  """
$ws.Range("E46").Copy()
$ws.Range("E46").PasteSpecial(-4163)
$ws.Range("E47").Formula = "=""  +0.41%  """
$ws.Range("E47").Copy()
$ws.Range("E47").PasteSpecial(-4163)
$ws.Range("D48").Formula = "=""1.941.47"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Formula = "=""  +0.04%  """
$ws.Range("E48").Copy()
$ws.Range("E48").PasteSpecial(-4163)
$ws.Range("E49").Formula = "=""  -1.17%  """
$ws.Range("E49").Copy()
$ws.Range("E49").PasteSpecial(-4163)
$ws.Range("E50").Formula = "=""  -3.56%  """
$ws.Range("E50").Copy()
$ws.Range("E50").PasteSpecial(-4163)
$ws.Range("D51").Formula = "=""114.15"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Formula = "=""  +1.92%  """
$ws.Range("E51").Copy()
$ws.Range("E51").PasteSpecial(-4163)
$excel.CutCopyMode = 0
